$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 5
$ws.Range("H2").Value = 5

# Row 9
$ws.Range("E9").Value = 10

# Row 14
$ws.Range("E14").Value = 20
$ws.Range("F14").Value = 7
$ws.Range("H14").Value = 7

# Row 16
$ws.Range("E16").Value = 184

# Row 17
$ws.Range("E17").Value = 9

# Row 18
$ws.Range("E18").Value = 42
